$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (row 23) down into
# the two new rows so the new cells pick up the same number-format styles
# (date style for column A, currency style for C:E, integer style for F)
# instead of creating brand-new style entries.
$ws.Range("A23:F23").Copy()
$ws.Range("A24:F24").PasteSpecial()
$ws.Range("A23:F23").Copy()
$ws.Range("A25:F25").PasteSpecial()

# Row 24 - 四方坪站, 2025-10-12
$ws.Cells.Item(24, 1).Value = 45942
$ws.Cells.Item(24, 2).Value = "四方坪站"
$ws.Cells.Item(24, 3).Value = 9225.0300000000007
$ws.Cells.Item(24, 4).Value = 7665.05
$ws.Cells.Item(24, 5).Value = 3175.58
$ws.Cells.Item(24, 6).Value = 375

# Row 25 - 高岭站, 2025-10-12
$ws.Cells.Item(25, 1).Value = 45942
$ws.Cells.Item(25, 2).Value = "高岭站"
$ws.Cells.Item(25, 3).Value = 3777.63
$ws.Cells.Item(25, 4).Value = 3006.63
$ws.Cells.Item(25, 5).Value = 991.63
$ws.Cells.Item(25, 6).Value = 134

# Reposition the view to match the author's final window state: scrolled so
# row 13 is at the top, with I23 as the active/selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("I23").Select()
